$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as a genuine text/string cell (avoids Excel
# auto-converting numeric-looking strings like "1.003" into numbers,
# and does so without leaving any NumberFormat/style residue behind).
function Set-TextValue($cellRef, $val) {
    $helper = $ws.Range("ZZ1")
    $escaped = $val.Replace('"', '""')
    $helper.Formula = '="' + $escaped + '"'
    $helper.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
    $helper.ClearContents()
}

$ws.Range('D2').Value = '27.686.50'
$ws.Range('E2').Value = '  +1.14%  '

$ws.Range('D3').Value = '1.874.16'
$ws.Range('E3').Value = '  +0.84%  '

Set-TextValue 'D4' '1.003'
$ws.Range('E4').Value = '  +0.17%  '

Set-TextValue 'D5' '331.27'
$ws.Range('E5').Value = '  +2.36%  '

Set-TextValue 'D6' '1.003'
$ws.Range('E6').Value = '  +0.16%  '

Set-TextValue 'D7' '0.4722'
$ws.Range('E7').Value = '  +4.15%  '

Set-TextValue 'D8' '0.3940'
$ws.Range('E8').Value = '  +1.72%  '

Set-TextValue 'D9' '47.92'

Set-TextValue 'D10' '0.08053'
$ws.Range('E10').Value = '  +1.47%  '

Set-TextValue 'D11' '1.028'
$ws.Range('E11').Value = '  +1.09%  '

Set-TextValue 'D12' '22.04'
$ws.Range('E12').Value = '  +2.90%  '

$ws.Range('D13').Value = '1.867.00'
$ws.Range('E13').Value = '  +0.37%  '

Set-TextValue 'D14' '5.955'
$ws.Range('E14').Value = '  +0.49%  '

Set-TextValue 'D15' '7.118'
$ws.Range('E15').Value = '  -0.24%  '

Set-TextValue 'D16' '1.006'
$ws.Range('E16').Value = '  +0.32%  '

$ws.Range('E17').Value = '  +1.61%  '

$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue 'D18' '0.06684'
$ws.Range('E18').Value = '  +2.52%  '

$ws.Range('B19').Value = 'Litecoin'
$ws.Range('C19').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue 'D19' '86.88'
$ws.Range('E19').Value = '  +0.93%  '

$ws.Range('E20').Value = '  +0.45%  '

$ws.Range('E21').Value = '  +0.20%  '

$ws.Range('D22').Value = '27.698.77'
$ws.Range('E22').Value = '  +1.17%  '

Set-TextValue 'D23' '5.512'
$ws.Range('E23').Value = '  -0.53%  '

$ws.Range('E24').Value = '  +0.86%  '

Set-TextValue 'D25' '2.305'
$ws.Range('E25').Value = '  +1.07%  '

$ws.Range('D26').Value = '2.111.06'
$ws.Range('E26').Value = '  +1.26%  '

Set-TextValue 'D27' '158.93'
$ws.Range('E27').Value = '  +3.26%  '

Set-TextValue 'D28' '20.10'
$ws.Range('E28').Value = '  +0.83%  '

Set-TextValue 'D29' '2.102'
$ws.Range('E29').Value = '  +0.96%  '

Set-TextValue 'D30' '5.567'
$ws.Range('E30').Value = '  +2.36%  '

Set-TextValue 'D31' '122.33'
$ws.Range('E31').Value = '  +0.88%  '

Set-TextValue 'D32' '0.9744'
$ws.Range('E32').Value = '  +3.89%  '

Set-TextValue 'D33' '0.09522'
$ws.Range('E33').Value = '  +2.42%  '

Set-TextValue 'D34' '1.448'
$ws.Range('E34').Value = '  -2.64%  '

Set-TextValue 'D35' '3.593'
$ws.Range('E35').Value = '  -0.12%  '

Set-TextValue 'D36' '5.325'
$ws.Range('E36').Value = '  +1.06%  '

Set-TextValue 'D37' '0.06103'
$ws.Range('E37').Value = '  +1.73%  '

$ws.Range('E38').Value = '  +0.48%  '

Set-TextValue 'D39' '1.226'
$ws.Range('E39').Value = '  -0.38%  '

Set-TextValue 'D40' '8.118'
$ws.Range('E40').Value = '  -0.97%  '

Set-TextValue 'D41' '0.6017'
$ws.Range('E41').Value = '  +1.58%  '

Set-TextValue 'D42' '0.1900'
$ws.Range('E42').Value = '  -0.15%  '

$ws.Range('E43').Value = '  +0.95%  '

Set-TextValue 'D44' '1.265'
$ws.Range('E44').Value = '  -1.28%  '

Set-TextValue 'D45' '0.5692'
$ws.Range('E45').Value = '  +1.26%  '

Set-TextValue 'D46' '12.21'
$ws.Range('E46').Value = '  +1.72%  '

Set-TextValue 'D47' '1.944'
$ws.Range('E47').Value = '  +0.92%  '

Set-TextValue 'D48' '3.379'
$ws.Range('E48').Value = '  +0.00%  '

Set-TextValue 'D49' '0.06893'
$ws.Range('E49').Value = '  +1.87%  '

Set-TextValue 'D50' '114.70'
$ws.Range('E50').Value = '  +5.79%  '

$ws.Range('B51').Value = 'EOS'
$ws.Range('C51').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
Set-TextValue 'D51' '1.070'
$ws.Range('E51').Value = '  +1.47%  '

$ws.Range("ZZ1").Clear() | Out-Null
